$d = $word.ActiveDocument

$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)   # wdHeaderFooterFirstPage
$shp = $hdr.Shapes.Item(3)
$tf = $shp.TextFrame
$tr = $tf.TextRange
for ($i = 1; $i -le $tr.Paragraphs.Count; $i++) {
    $p = $tr.Paragraphs.Item($i)
    Write-Host "Para $i :" $p.Range.Text
}
